$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the bottom four data rows (rows 4-7) - only 3 students remain
# ------------------------------------------------------------------
$ws.Rows.Item(4).Resize(4).Delete()

# ------------------------------------------------------------------
# 2. Remove column D ("2025-03-06 - hw4") - column E shifts to D,
#    column F shifts to E
# ------------------------------------------------------------------
$ws.Columns.Item(4).Delete()

# ------------------------------------------------------------------
# 3. Update header text in (new) column D
# ------------------------------------------------------------------
$ws.Range("D1").Value = "2025-03-13 - asdasdadaasd"

# ------------------------------------------------------------------
# 4. Recode attendance marks to numeric flags (0 = absent, 1 = present)
#    Values must stay text (matches the sheet's numberStoredAsText data),
#    so force a text format, assign, then drop the format again so no
#    stray style survives on the cell.
# ------------------------------------------------------------------
$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "0"
$d2.ClearFormats()

$e2 = $ws.Range("E2")
$e2.NumberFormat = "@"
$e2.Value = "1"
$e2.ClearFormats()

$d3 = $ws.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "0"
$d3.ClearFormats()

$e3 = $ws.Range("E3")
$e3.NumberFormat = "@"
$e3.Value = "0"
$e3.ClearFormats()

# ------------------------------------------------------------------
# 5. Resize columns to match the new layout
#    (column 1 already has the right width from the source file and is
#    intentionally left untouched - re-assigning it would push it through
#    Excel's pixel-snapping and change its stored value)
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 27
$ws.Columns.Item(5).ColumnWidth = 18

Write-Output "edit complete"
